# Daily attendance processing - 2025-11-27 20:26:23
# Normalizes the "Recorded By" (column G) entries: when the comma-separated
# list of recorders ends with a "System"/"system" entry, the whole list
# order is reversed so the System-like entry moves to the front.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $raw = $cell.Value2

    if ($raw -eq $null) {
        continue
    }

    $text = [string]$raw
    if ($text -eq "") {
        continue
    }

    $parts = $text -split ","
    for ($i = 0; $i -lt $parts.Length; $i++) {
        $parts[$i] = $parts[$i].Trim()
    }

    $lastVal = $parts[$parts.Length - 1]

    if ($parts.Length -ge 2 -and $lastVal.ToLower() -eq "system") {
        $reversed = $parts[($parts.Length - 1)..0]
        $newText = [string]::Join(", ", $reversed)
        $cell.Value = $newText
    }
}
